$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 currently holds the "old" record (D27=44323 date, M27=270 volume).
# The edit keeps that old record but moves it to row 28, and row 27 gets
# the updated values (D27=45127, M27=200). All other columns stay the same
# between the two rows.

# Copy each cell value (and number format) from row 27 into row 28.
for ($col = 1; $col -le 20; $col++) {
    $srcCell = $ws.Cells.Item(27, $col)
    $dstCell = $ws.Cells.Item(28, $col)
    $dstCell.Value = $srcCell.Value2
}
# Column D (4) is a date formatted cell on row 27; match that format on row 28.
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat

# Now update row 27 with the new values
$ws.Range("D27").Value = 45127
$ws.Range("M27").Value = 200
